$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format for Price column (D) so numeric-looking strings
# (e.g. "1.00", "36.77") are preserved exactly as text, matching the
# original inline-string cell values instead of being parsed as numbers.
$priceTextCells = @(
    4, 5, 6, 7, 9, 10, 11, 12, 13, 15, 17, 20, 21, 24, 25, 26, 27, 28, 29, 30, 31, 32, 33, 34, 36, 39, 40, 41, 42, 43, 45, 46, 48, 49, 50, 51
)
foreach ($r in $priceTextCells) {
    $ws.Cells.Item($r, 4).NumberFormat = "@"
}

# Row 2
$ws.Cells.Item(2, 4).Value = "42.659.62"
$ws.Cells.Item(2, 5).Value = "  -6.80%  "

# Row 3
$ws.Cells.Item(3, 4).Value = "2.220.22"
$ws.Cells.Item(3, 5).Value = "  -7.48%  "

# Row 4
$ws.Cells.Item(4, 4).Value = "1.00"
$ws.Cells.Item(4, 5).Value = "  +0.05%  "

# Row 5
$ws.Cells.Item(5, 4).Value = "312.76"
$ws.Cells.Item(5, 5).Value = "  -1.78%  "

# Row 6
$ws.Cells.Item(6, 4).Value = "98.30"
$ws.Cells.Item(6, 5).Value = "  -13.45%  "

# Row 7
$ws.Cells.Item(7, 4).Value = "0.571"
$ws.Cells.Item(7, 5).Value = "  -9.93%  "

# Row 8
$ws.Cells.Item(8, 5).Value = "  +0.05%  "

# Row 9
$ws.Cells.Item(9, 4).Value = "0.558"
$ws.Cells.Item(9, 5).Value = "  -11.03%  "

# Row 10
$ws.Cells.Item(10, 4).Value = "36.77"
$ws.Cells.Item(10, 5).Value = "  -12.15%  "

# Row 11
$ws.Cells.Item(11, 4).Value = "53.83"
$ws.Cells.Item(11, 5).Value = "  -4.12%  "

# Row 12
$ws.Cells.Item(12, 4).Value = "0.0835"
$ws.Cells.Item(12, 5).Value = "  -10.03%  "

# Row 13
$ws.Cells.Item(13, 4).Value = "7.59"
$ws.Cells.Item(13, 5).Value = "  -13.06%  "

# Row 14
$ws.Cells.Item(14, 5).Value = "  -4.74%  "

# Row 15
$ws.Cells.Item(15, 4).Value = "0.869"
$ws.Cells.Item(15, 5).Value = "  -13.50%  "

# Row 16
$ws.Cells.Item(16, 4).Value = "2.554.24"
$ws.Cells.Item(16, 5).Value = "  -7.46%  "

# Row 17
$ws.Cells.Item(17, 4).Value = "13.92"
$ws.Cells.Item(17, 5).Value = "  -11.97%  "

# Row 18
$ws.Cells.Item(18, 4).Value = "2.220.66"
$ws.Cells.Item(18, 5).Value = "  -7.19%  "

# Row 19
$ws.Cells.Item(19, 4).Value = "42.539.72"
$ws.Cells.Item(19, 5).Value = "  -6.86%  "

# Row 20
$ws.Cells.Item(20, 4).Value = "13.82"
$ws.Cells.Item(20, 5).Value = "  +3.14%  "

# Row 21
$ws.Cells.Item(21, 4).Value = "6.59"
$ws.Cells.Item(21, 5).Value = "  -11.75%  "

# Row 22
$ws.Cells.Item(22, 4).Value = "0.0₃0946"
$ws.Cells.Item(22, 5).Value = "  -12.72%  "

# Row 24
$ws.Cells.Item(24, 4).Value = "64.63"
$ws.Cells.Item(24, 5).Value = "  -13.51%  "

# Row 25
$ws.Cells.Item(25, 4).Value = "234.05"
$ws.Cells.Item(25, 5).Value = "  -11.75%  "

# Row 26
$ws.Cells.Item(26, 4).Value = "2.11"
$ws.Cells.Item(26, 5).Value = "  -10.11%  "

# Row 27
$ws.Cells.Item(27, 4).Value = "1.00"
$ws.Cells.Item(27, 5).Value = "  +0.05%  "

# Row 28
$ws.Cells.Item(28, 4).Value = "10.15"
$ws.Cells.Item(28, 5).Value = "  -10.15%  "

# Row 29
$ws.Cells.Item(29, 2).Value = "Filecoin"
$ws.Cells.Item(29, 3).Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Cells.Item(29, 4).Value = "6.50"
$ws.Cells.Item(29, 5).Value = "  -14.31%  "

# Row 30
$ws.Cells.Item(30, 2).Value = "Toncoin"
$ws.Cells.Item(30, 3).Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Cells.Item(30, 4).Value = "2.15"
$ws.Cells.Item(30, 5).Value = "  -8.65%  "

# Row 31
$ws.Cells.Item(31, 4).Value = "0.0880"
$ws.Cells.Item(31, 5).Value = "  -9.55%  "

# Row 32
$ws.Cells.Item(32, 4).Value = "20.41"
$ws.Cells.Item(32, 5).Value = "  -10.34%  "

# Row 33
$ws.Cells.Item(33, 4).Value = "156.77"
$ws.Cells.Item(33, 5).Value = "  -9.24%  "

# Row 34
$ws.Cells.Item(34, 4).Value = "32.96"
$ws.Cells.Item(34, 5).Value = "  -15.41%  "

# Row 35
$ws.Cells.Item(35, 5).Value = "  -9.40%  "

# Row 36
$ws.Cells.Item(36, 4).Value = "3.06"
$ws.Cells.Item(36, 5).Value = "  +1.30%  "

# Row 37
$ws.Cells.Item(37, 5).Value = "  -7.65%  "

# Row 38
$ws.Cells.Item(38, 5).Value = "  -10.19%  "

# Row 39
$ws.Cells.Item(39, 2).Value = "Kaspa"
$ws.Cells.Item(39, 3).Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Cells.Item(39, 4).Value = "0.105"
$ws.Cells.Item(39, 5).Value = "  -10.85%  "

# Row 40
$ws.Cells.Item(40, 2).Value = "ARBITRUM"
$ws.Cells.Item(40, 3).Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Cells.Item(40, 4).Value = "1.82"
$ws.Cells.Item(40, 5).Value = "  +2.16%  "

# Row 41
$ws.Cells.Item(41, 2).Value = "VeChain"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Cells.Item(41, 4).Value = "0.0320"
$ws.Cells.Item(41, 5).Value = "  -11.60%  "

# Row 42
$ws.Cells.Item(42, 2).Value = "NEARProtocol"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Cells.Item(42, 4).Value = "3.49"
$ws.Cells.Item(42, 5).Value = "  -15.82%  "

# Row 43
$ws.Cells.Item(43, 4).Value = "1.00"
$ws.Cells.Item(43, 5).Value = "  +0.24%  "

# Row 44
$ws.Cells.Item(44, 4).Value = "1.784.31"
$ws.Cells.Item(44, 5).Value = "  +7.33%  "

# Row 45
$ws.Cells.Item(45, 4).Value = "88.69"
$ws.Cells.Item(45, 5).Value = "  -14.22%  "

# Row 46
$ws.Cells.Item(46, 4).Value = "11.98"
$ws.Cells.Item(46, 5).Value = "  -11.50%  "

# Row 47
$ws.Cells.Item(47, 5).Value = "  -14.82%  "

# Row 48
$ws.Cells.Item(48, 4).Value = "77.22"
$ws.Cells.Item(48, 5).Value = "  -11.69%  "

# Row 49
$ws.Cells.Item(49, 4).Value = "5.37"
$ws.Cells.Item(49, 5).Value = "  -4.95%  "

# Row 50
$ws.Cells.Item(50, 4).Value = "60.07"
$ws.Cells.Item(50, 5).Value = "  -16.39%  "

# Row 51
$ws.Cells.Item(51, 4).Value = "8.52"
$ws.Cells.Item(51, 5).Value = "  -9.86%  "
